# Updated test data for 5,24,40V,BatteryStandby and AC Calculations test cases
#
# Sheet "Add Panels and Devices":
#   - B4  (User Story ticket id): NGC-1679 -> NGC-571/T1402 OR TC-186
#        (also drops the highlighted/font style back to the default style)
#   - A8  (Panel Name):            MX 4000 -> MZX252
#   - C8  (CPU Type):              CPU 800 -> "" (cleared, kept as a text cell)
#   - F8  (Expected 24V PSU Load): 0.205 (number) -> "0.100" (entered as text,
#        i.e. with a leading apostrophe / quote-prefix, same as the authored file)
#   - J8, N8, O8 (Expected 2nd/3rd/4th 24V PSU Load values): 0.207 -> 0.165
#   - Active selection moves from P8 to B8

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# User Story / ticket reference
$ws.Range("B4").Value = "NGC-571/T1402 OR TC-186"
$ws.Range("B4").Style = "Normal"

# Row 8 test data
$ws.Range("A8").Value = "MZX252"
$ws.Range("C8").Value = "'"
$ws.Range("F8").Value = "'0.100"
$ws.Range("J8").Value = 0.165
$ws.Range("N8").Value = 0.165
$ws.Range("O8").Value = 0.165

# Match the saved selection/active cell
[void]$ws.Range("B8").Select()
